$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.415.58"
$ws.Range("E2").Value = "  +3.69%  "

$ws.Range("D3").Value = "1.838.87"
$ws.Range("E3").Value = "  +3.65%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.025"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +1.92%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "318.59"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +4.01%  "

$ws.Range("E6").Value = "  +2.10%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4358"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +1.98%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3726"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +2.76%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07336"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +2.70%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.8723"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +3.94%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "21.32"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +3.83%  "

$ws.Range("D12").Value = "1.949.96"
$ws.Range("E12").Value = "  +7.86%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.471"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +4.39%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.677"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +3.64%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.07134"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +3.48%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "82.16"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +4.37%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.030"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +2.13%  "

$ws.Range("E18").Value = "  +3.40%  "

$ws.Range("E19").Value = "  +2.33%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "15.39"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +3.08%  "

$ws.Range("D21").Value = "27.439.39"
$ws.Range("E21").Value = "  +3.73%  "

$ws.Range("E22").Value = "  +2.83%  "

$ws.Range("E23").Value = "  +1.04%  "

$ws.Range("D24").Value = "2.149.98"
$ws.Range("E24").Value = "  +6.24%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "156.87"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +2.88%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.894"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +2.00%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "18.55"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +3.14%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "5.244"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +3.87%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.915"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +7.99%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "115.31"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.29%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.09046"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.77%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.198"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +7.52%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.7576"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +4.73%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.461"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +3.60%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.860"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +4.10%  "

$ws.Range("E36").Value = "  +2.31%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.149"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +4.47%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01954"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +3.63%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.05239"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +2.19%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.5168"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +5.43%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.779"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +6.47%  "

$ws.Range("E42").Value = "  +3.14%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "6.530"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +3.22%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "8.449"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +6.18%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "108.43"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +3.63%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "10.45"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +3.25%  "

$ws.Range("E47").Value = "  +2.42%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.4628"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +3.67%  "

$ws.Range("E49").Value = "  +2.47%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.06290"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.74%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.867"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +9.27%  "
